# Apply the "Updated cryptos list" data refresh to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D (Price) while forcing it to stay a
# plain text cell (matches the source data, which is inline/shared string,
# never a real number) and without leaving a lingering custom cell style.
function Set-PriceText([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText "D2" "60.940.93"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3 - Ethereum
Set-PriceText "D3" "3.389.78"
$ws.Range("E3").Value = "  -1.28%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-PriceText "D5" "570.95"
$ws.Range("E5").Value = "  -0.79%  "

# Row 6 - Solana
Set-PriceText "D6" "142.01"
$ws.Range("E6").Value = "  -2.18%  "

# Row 7 - was USDC, now LidoStakedEther
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-PriceText "D7" "3.389.58"
$ws.Range("E7").Value = "  -1.35%  "

# Row 8 - was LidoStakedEther, now USDC
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-PriceText "D8" "1.00"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.48%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -1.43%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.50%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +2.18%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-PriceText "D13" "3.968.30"
$ws.Range("E13").Value = "  -1.23%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.92%  "

# Row 15 - Avalanche
Set-PriceText "D15" "28.28"
$ws.Range("E15").Value = "  +1.15%  "

# Row 16 - was WrappedEther, now ShibaInu
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-PriceText "D16" "0.0000171"
$ws.Range("E16").Value = "  -1.02%  "

# Row 17 - was ShibaInu, now WrappedEther
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-PriceText "D17" "3.393.01"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18 - WrappedBTC
Set-PriceText "D18" "60.989.02"

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -1.83%  "

# Row 20 - Chainlink
Set-PriceText "D20" "13.88"
$ws.Range("E20").Value = "  -2.20%  "

# Row 21 - Uniswap
Set-PriceText "D21" "8.99"
$ws.Range("E21").Value = "  -4.62%  "

# Row 22 - BitcoinCash
Set-PriceText "D22" "384.53"
$ws.Range("E22").Value = "  -2.73%  "

# Row 23 - Polygon
Set-PriceText "D23" "0.557"
$ws.Range("E23").Value = "  -1.32%  "

# Row 24 - Litecoin
Set-PriceText "D24" "74.35"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.59%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -4.16%  "

# Row 27 - WrappedeETH
Set-PriceText "D27" "3.527.75"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -0.94%  "

# Row 30 - RenderToken
$ws.Range("E30").Value = "  -2.83%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("E31").Value = "  -2.77%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  -2.19%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  -1.68%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  +0.03%  "

# Row 35 - EthereumClassic
Set-PriceText "D35" "23.51"
$ws.Range("E35").Value = "  -1.56%  "

# Row 36 - Aptos
$ws.Range("E36").Value = "  -0.70%  "

# Row 37 - Monero
Set-PriceText "D37" "167.27"
$ws.Range("E37").Value = "  -0.16%  "

# Row 38 - RenzoRestakedETH
Set-PriceText "D38" "3.419.84"
$ws.Range("E38").Value = "  -1.17%  "

# Row 39 - NEARProtocol
$ws.Range("E39").Value = "  -2.48%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  -4.60%  "

# Row 41 - Hedera
Set-PriceText "D41" "0.0775"
$ws.Range("E41").Value = "  -1.27%  "

# Row 42 - EnergySwap
Set-PriceText "D42" "27.31"
$ws.Range("E42").Value = "  +1.98%  "

# Row 43 - was FirstDigitalUSD, now Mantle
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceText "D43" "0.780"
$ws.Range("E43").Value = "  -2.31%  "

# Row 44 - was Mantle, now FirstDigitalUSD
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-PriceText "D44" "0.999"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45 - OKB
$ws.Range("E45").Value = "  -0.27%  "

# Row 46 - Filecoin
Set-PriceText "D46" "4.43"
$ws.Range("E46").Value = "  -1.37%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  -3.37%  "

# Row 48 - ONDO
Set-PriceText "D48" "1.14"
$ws.Range("E48").Value = "  -1.42%  "

# Row 49 - Maker
Set-PriceText "D49" "2.484.69"
$ws.Range("E49").Value = "  -4.19%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  -1.44%  "

# Row 51 - InjectiveProtocol
Set-PriceText "D51" "22.99"
$ws.Range("E51").Value = "  -0.69%  "
